$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D2:D51) keeps its original text representation,
# since several values (e.g. "30.827.37", "1.190", "1.000") must remain as
# literal strings rather than being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.827.37'
$ws.Range("E2").Value = '  +2.15%  '
$ws.Range("D3").Value = '2.120.40'
$ws.Range("E3").Value = '  +10.26%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '334.57'
$ws.Range("E5").Value = '  +4.67%  '
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = '0.5271'
$ws.Range("E7").Value = '  +4.01%  '
$ws.Range("E8").Value = '  +8.32%  '
$ws.Range("D9").Value = '0.09078'
$ws.Range("E9").Value = '  +8.55%  '
$ws.Range("D10").Value = '46.91'
$ws.Range("E10").Value = '  +10.76%  '
$ws.Range("D11").Value = '1.190'
$ws.Range("E11").Value = '  +6.90%  '
$ws.Range("D12").Value = '25.41'
$ws.Range("E12").Value = '  +5.64%  '
$ws.Range("D13").Value = '2.112.89'
$ws.Range("E13").Value = '  +10.07%  '
$ws.Range("D14").Value = '6.776'
$ws.Range("E14").Value = '  +5.44%  '
$ws.Range("D15").Value = '7.859'
$ws.Range("E15").Value = '  +8.17%  '
$ws.Range("D16").Value = '98.31'
$ws.Range("E16").Value = '  +6.08%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001142'
$ws.Range("E17").Value = '  +4.14%  '
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '0.06639'
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '19.23'
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("D21").Value = '6.411'
$ws.Range("E21").Value = '  +7.60%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = '30.932.93'
$ws.Range("E23").Value = '  +2.52%  '
$ws.Range("D24").Value = '12.18'
$ws.Range("E24").Value = '  +7.08%  '
$ws.Range("D25").Value = '2.366.13'
$ws.Range("E25").Value = '  +10.72%  '
$ws.Range("D26").Value = '2.264'
$ws.Range("E26").Value = '  +3.15%  '
$ws.Range("D27").Value = '22.99'
$ws.Range("E27").Value = '  +4.68%  '
$ws.Range("D28").Value = '2.589'
$ws.Range("E28").Value = '  +14.06%  '
$ws.Range("D29").Value = '163.61'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("D30").Value = '133.70'
$ws.Range("E30").Value = '  +3.60%  '
$ws.Range("D31").Value = '1.182'
$ws.Range("E31").Value = '  +4.47%  '
$ws.Range("D32").Value = '0.1078'
$ws.Range("E32").Value = '  +2.89%  '
$ws.Range("D33").Value = '6.267'
$ws.Range("E33").Value = '  +4.96%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '4.019'
$ws.Range("E34").Value = '  +5.89%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.567'
$ws.Range("E35").Value = '  +30.02%  '
$ws.Range("D36").Value = '0.02603'
$ws.Range("E36").Value = '  +5.95%  '
$ws.Range("E37").Value = '  +4.99%  '
$ws.Range("D38").Value = '9.631'
$ws.Range("E38").Value = '  +12.13%  '
$ws.Range("D39").Value = '0.06753'
$ws.Range("E39").Value = '  +4.82%  '
$ws.Range("D40").Value = '12.74'
$ws.Range("E40").Value = '  +10.99%  '
$ws.Range("D41").Value = '0.2273'
$ws.Range("E41").Value = '  +5.66%  '
$ws.Range("D42").Value = '0.6866'
$ws.Range("E42").Value = '  +5.22%  '
$ws.Range("D43").Value = '1.262'
$ws.Range("E43").Value = '  +4.05%  '
$ws.Range("D44").Value = '14.22'
$ws.Range("E44").Value = '  +6.11%  '
$ws.Range("D45").Value = '0.6476'
$ws.Range("E45").Value = '  +6.94%  '
$ws.Range("D46").Value = '0.9988'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").Value = '2.264'
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").Value = '1.289'
$ws.Range("E49").Value = '  +6.28%  '
$ws.Range("D50").Value = '83.21'
$ws.Range("E50").Value = '  +5.24%  '
$ws.Range("D51").Value = '0.07087'
$ws.Range("E51").Value = '  +3.69%  '
